# Insert a new data row at spreadsheet row 69 (pushing the existing row 69
# and all rows below it down by one, row 184 becomes row 185), then fill
# in the values for the newly inserted record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(69).Insert()

$ws.Range("A69").Value2 = 4
$ws.Range("B69").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C69").Value2 = "Los Lagos"
$ws.Range("D69").Value2 = 44533
$ws.Range("E69").Value2 = 10
$ws.Range("F69").Value2 = 100112044
$ws.Range("G69").Value2 = "Perejil"
$ws.Range("H69").Value2 = "Sin especificar"
$ws.Range("I69").Value2 = "Primera"
$ws.Range("J69").Value2 = 180
$ws.Range("K69").Value2 = 5000
$ws.Range("L69").Value2 = 5000
$ws.Range("M69").Value2 = 5000
$ws.Range("N69").Value2 = "`$/docena de atados (3 kilos)"
$ws.Range("O69").Value2 = "Región Metropolitana"
$ws.Range("P69").Value2 = 1667
$ws.Range("Q69").Value2 = 3
$ws.Range("R69").Value2 = "Hortaliza"
